# Apply cryptos list update (Thu Jun 20 17:45:47 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.745.94'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.50%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.517.20'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.27%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.76'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.10'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.65%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.516.97'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.30%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  -0.95%  '
$ws.Range("E10").Value = '  +0.64%  '
$ws.Range("E11").Value = '  +2.49%  '
$ws.Range("E12").Value = '  -0.88%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.119.12'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.13%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.60'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.90%  '
$ws.Range("E15").Value = '  -0.72%  '
$ws.Range("E16").Value = '  +0.23%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.516.10'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.814.57'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.04'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.37'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("E21").Value = '  -2.98%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '391.88'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.660.92'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.04'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.24%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000112'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -3.71%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.65'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.76%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.58'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +10.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.13%  '
$ws.Range("E31").Value = '  -0.71%  '
$ws.Range("E32").Value = '  -0.14%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.523.63'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '24.18'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.54%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.25'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +4.62%  '
$ws.Range("E38").Value = '  +0.94%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.91'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '168.60'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.85%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0817'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.69%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.821'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.84%  '
$ws.Range("B43").Value = 'ONDO'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.25'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.91%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '25.61'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -3.16%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.67'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("E47").Value = '  -0.98%  '
$ws.Range("E48").Value = '  -1.97%  '
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.385.78'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -4.37%  '
$ws.Range("E51").Value = '  +0.94%  '
